$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-02-11 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-02-12 Monday", 2) | Out-Null

# Update the answer table cells (positional, to avoid collisions between old/new values)
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "85÷3=28, 1"  # was "56÷9=6, 2"
$t.Cell(1, 2).Range.Text = "25÷3=8, 1"  # was "20÷4=5, 0"
$t.Cell(1, 3).Range.Text = "97÷2=48, 1"  # was "81÷8=10, 1"
$t.Cell(1, 4).Range.Text = "15÷7=2, 1"  # was "46÷5=9, 1"
$t.Cell(1, 5).Range.Text = "16÷9=1, 7"  # was "47÷5=9, 2"
$t.Cell(5, 1).Range.Text = "77÷4=19, 1"  # was "18÷2=9, 0"
$t.Cell(5, 2).Range.Text = "20÷3=6, 2"  # was "35÷2=17, 1"
$t.Cell(5, 3).Range.Text = "76÷4=19, 0"  # was "32÷2=16, 0"
$t.Cell(5, 4).Range.Text = "76÷2=38, 0"  # was "31÷8=3, 7"
$t.Cell(5, 5).Range.Text = "96÷3=32, 0"  # was "85÷5=17, 0"
$t.Cell(9, 1).Range.Text = "44÷2=22, 0"  # was "95÷6=15, 5"
$t.Cell(9, 2).Range.Text = "35÷2=17, 1"  # was "48÷6=8, 0"
$t.Cell(9, 3).Range.Text = "82÷7=11, 5"  # was "63÷9=7, 0"
$t.Cell(9, 4).Range.Text = "27÷5=5, 2"  # was "97÷2=48, 1"
$t.Cell(9, 5).Range.Text = "34÷5=6, 4"  # was "52÷3=17, 1"
$t.Cell(13, 1).Range.Text = "59÷6=9, 5"  # was "13÷2=6, 1"
$t.Cell(13, 2).Range.Text = "23÷9=2, 5"  # was "95÷7=13, 4"
$t.Cell(13, 3).Range.Text = "16÷3=5, 1"  # was "48÷7=6, 6"
$t.Cell(13, 4).Range.Text = "34÷7=4, 6"  # was "86÷6=14, 2"
$t.Cell(13, 5).Range.Text = "82÷5=16, 2"  # was "60÷6=10, 0"
$t.Cell(17, 1).Range.Text = "79÷3=26, 1"  # was "40÷8=5, 0"
$t.Cell(17, 2).Range.Text = "16÷5=3, 1"  # was "16÷9=1, 7"
$t.Cell(17, 3).Range.Text = "57÷4=14, 1"  # was "18÷4=4, 2"
$t.Cell(17, 4).Range.Text = "65÷5=13, 0"  # was "57÷3=19, 0"
$t.Cell(17, 5).Range.Text = "56÷5=11, 1"  # was "22÷8=2, 6"
